# PowerShell-style Word COM-interop script applying the edits described
# by the unified diff (Redovisning rapport.docx).
#
# Most hunks in the diff are just proofErr (spell-check squiggle) markers
# disappearing and adjacent runs merging -- that happens naturally in Word
# whenever the covered text is retyped/edited, so we reproduce it here by
# doing ordinary Find&Replace (wdFindContinue) over the same text; Word's
# replace collapses the old runs (and their proofErr wrappers) into a
# fresh run. Only the replacements below actually change the rendered
# text.

$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2) | Out-Null
}

# "Anna Ng, Cam H ..." -> "Anna N, Cam H ..."
Replace-Text "Anna Ng, Cam H" "Anna N, Cam H"

# "Cam: programmering och fördelning ..." (proofErr removal only, re-set
# text so the run/proofErr collapse like the other hunks)
Replace-Text "Cam: programmering och fördelning av arbetsuppgift till andra medlemmar" "Cam: programmering och fördelning av arbetsuppgift till andra medlemmar"

# "Cam:  Deltagit i alla gruppmöte" (proofErr removal only)
Replace-Text "Cam:  Deltagit i alla gruppmöte" "Cam:  Deltagit i alla gruppmöte"

# Anna's "Uppdelning och uppföljning" paragraph (proofErr removal only)
Replace-Text "Var ansvarig för att ta fram grafiska bilder på alla klasser och modeller, var även ansvarig för kodning av validering funktionerna för GUI, design av Customer och Transaction klassen samt JUnit test av klassen." "Var ansvarig för att ta fram grafiska bilder på alla klasser och modeller, var även ansvarig för kodning av validering funktionerna för GUI, design av Customer och Transaction klassen samt JUnit test av klassen."

# Cam's "Uppdelning och uppföljning" paragraph (proofErr removal only)
Replace-Text "Cam: Var ansvarig för BankLogic och Repository klassen samt JUnit test av dessa klasser, skapade databasscripten, var delaktig i designen av kodningen till GUI samt testade applikationen. " "Cam: Var ansvarig för BankLogic och Repository klassen samt JUnit test av dessa klasser, skapade databasscripten, var delaktig i designen av kodningen till GUI samt testade applikationen. "

# Hampus's "Uppdelning och uppföljning" paragraph (proofErr removal only)
Replace-Text "Ansvarade för att generera Java dokumenten och JAR filen, testandet av applikationen samt delaktig i designen av kodning till GUI’n, designen av Customer klassen samt JUnit test av klasserna." "Ansvarade för att generera Java dokumenten och JAR filen, testandet av applikationen samt delaktig i designen av kodning till GUI’n, designen av Customer klassen samt JUnit test av klasserna."

# Åsa's "Uppdelning och uppföljning" paragraph: proofErr removal plus a
# new trailing sentence.
Replace-Text "Åsa: Ansvarig för design av klassen Account, SavingsAccount samt CreditAccount, skapade även testklasserna till klasserna. Testade applikationen, satte även upp konto för GitHub." "Åsa: Ansvarig för design av klassen Account, SavingsAccount samt CreditAccount, skapade även testklasserna till klasserna. Testade applikationen, satte även upp konto för GitHub. Redigerat i Javadoc."

# Hamid's "Uppdelning och uppföljning" paragraph (proofErr removal only)
Replace-Text ", ansvarig för SavingsAccount, CustomerTest, skötte resterande kodning som behövdes för att komplettera de större koderna, design av databas." ", ansvarig för SavingsAccount, CustomerTest, skötte resterande kodning som behövdes för att komplettera de större koderna, design av databas."

# "Arbetat enligt Scrum samtgjorde ..." -> insert missing space
Replace-Text "Arbetat enligt Scrum samtgjorde återkoppling" "Arbetat enligt Scrum samt gjorde återkoppling"

# "Planerade möten samt mail kontakt." -> drop the space ("mailkontakt")
Replace-Text "Planerade möten samt mail kontakt." "Planerade möten samt mailkontakt."

# "Kravshantering" heading -> "Kravhantering"
Replace-Text "Kravshantering" "Kravhantering"

# "på pdf dokumenten" -> "på pdf-dokumenten"
Replace-Text "på pdf dokumenten" "på pdf-dokumenten"

# "Krav och vad som behövs göra för att lösa" -> "... görs för att lösa"
Replace-Text "Krav och vad som behövs göra för att lösa" "Krav och vad som behövs göras för att lösa"

# "i Github ärendeshantering s.k." (proofErr removal only)
Replace-Text "i Github ärendeshantering s.k. ISSUE." "i Github ärendeshantering s.k. ISSUE."

# ", den som vill tar de hand om bugghanteringen." -> drop "de "
Replace-Text ", den som vill tar de hand om bugghanteringen." ", den som vill tar hand om bugghanteringen."

# "Versionhantering" heading (proofErr removal only)
Replace-Text "Versionhantering" "Versionhantering"

# "Vi höll oss uppdaterade ... GitHub." (proofErr removal only)
Replace-Text "Vi höll oss uppdaterade genom att pusha upp de nya koderna på GitHub." "Vi höll oss uppdaterade genom att pusha upp de nya koderna på GitHub."

# "Lessons learned" heading (proofErr removal only)
Replace-Text "Lessons learned" "Lessons learned"

# Final "Lessons learned" paragraph: several small wording tweaks.
Replace-Text "Alla hade inte samma Java kompetens" "Alla hade inte samma Javakompetens"
Replace-Text "så det var inte lika lätt att kunna tilldela" "så det var inte så lätt att kunna tilldela"
Replace-Text "vi kunde bli färdiga inom tid med så få problem som möjligt." "vi kunde bli färdiga i tid med så få problem som möjligt."
